{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Intent (per commit message \"Modified handling of empty SeaKen results\"\n// and the supplied OOXML diff): remove the stray w:proofErr\n// (spellStart/spellEnd/gramStart/gramEnd) markers that Word's proofing\n// pass had inserted around certain runs, and merge the runs that had\n// been artificially split by those markers back into single runs with\n// the same visible text. No visible text content changes.\n\nconst OOXML_NS =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>{BODY}</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapParagraph(innerParagraphXml) {\n  return OOXML_NS.replace('{BODY}', innerParagraphXml);\n}\n\n// Find the `occurrence`-th paragraph (1-based) whose trimmed text exactly\n// equals `text`, among the paragraphs already loaded in `items`.\nfunction findParagraphByText(items, text, occurrence) {\n  let seen = 0;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) {\n      seen++;\n      if (seen === occurrence) {\n        return items[i];\n      }\n    }\n  }\n  return null;\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) First table: \"Station_ID\" cell - drop spellStart/spellEnd around the run.\nconst stationId1 = findParagraphByText(items, \"Station_ID\", 1);\nif (stationId1) {\n  stationId1.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\">' +\n        \"<w:r><w:t>Station_ID</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 2) First table: \"Station_Description\" cell.\nconst stationDesc = findParagraphByText(items, \"Station_Description\", 1);\nif (stationDesc) {\n  stationDesc.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\">' +\n        \"<w:r><w:t>Station_Description</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 3) First table: \"SampleType\" cell.\nconst sampleType = findParagraphByText(items, \"SampleType\", 1);\nif (sampleType) {\n  sampleType.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\">' +\n        \"<w:r><w:t>SampleType</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 4) First table: \"StatusIdentifier\" cell.\nconst statusIdentifier = findParagraphByText(items, \"StatusIdentifier\", 1);\nif (statusIdentifier) {\n  statusIdentifier.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\">' +\n        \"<w:r><w:t>StatusIdentifier</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 5) \"Seasonal Kendall Results: ...\" paragraph - merge the 3 runs that\n// spelled out \"...seasonal \", \"kendall\", \" trend analysis...\" (split by\n// spellStart/gramStart/spellEnd/gramEnd around \"kendall\") into one run.\nconst seaKen = findParagraphByText(\n  items,\n  \"Seasonal Kendall Results: Table detailing results of seasonal kendall trend analysis for each parameter at each station. Column descriptions follow:\",\n  1\n);\nif (seaKen) {\n  seaKen.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\">' +\n        '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Seasonal Kendall Results: </w:t></w:r>' +\n        \"<w:r><w:t>Table detailing results of seasonal kendall trend analysis for each parameter at each station. Column descriptions follow:</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 6) Second table: \"Station_ID\" cell (2nd occurrence of that exact text).\nconst stationId2 = findParagraphByText(items, \"Station_ID\", 2);\nif (stationId2) {\n  stationId2.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\">' +\n        \"<w:r><w:t>Station_ID</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 7) Second table: \"Pvalue\" cell.\nconst pvalue = findParagraphByText(items, \"Pvalue\", 1);\nif (pvalue) {\n  pvalue.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\">' +\n        \"<w:r><w:t>Pvalue</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 8) Second table: \"signif\" cell.\nconst signif = findParagraphByText(items, \"signif\", 1);\nif (signif) {\n  signif.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\">' +\n        \"<w:r><w:t>signif</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 9) \"E. Coli/Enterococcus: \" heading - merge 3 runs into 1.\nconst eColi = findParagraphByText(items, \"E. Coli/Enterococcus: \", 1);\nif (eColi) {\n  eColi.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"00CC7C13\" w:rsidRDefault=\"00CC7C13\" w:rsidP=\"00CC7C13\">' +\n        '<w:pPr><w:pStyle w:val=\"Heading4\"/></w:pPr>' +\n        '<w:r><w:t xml:space=\"preserve\">E. Coli/Enterococcus: </w:t></w:r>' +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 10) \"pH:\" heading - merge 2 runs into 1 (keep lastRenderedPageBreak).\nconst phHeading = findParagraphByText(items, \"pH:\", 1);\nif (phHeading) {\n  phHeading.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"00CC7C13\" w:rsidRDefault=\"00CC7C13\" w:rsidP=\"00CC7C13\">' +\n        '<w:pPr><w:pStyle w:val=\"Heading4\"/></w:pPr>' +\n        \"<w:r><w:lastRenderedPageBreak/><w:t>pH:</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\n// 11) \"There are two options for plotting pH. The first is ...\" paragraph -\n// merge the first 3 runs (\"There are two options for plotting \", \"pH.\",\n// \" The first is\") into a single run; the trailing \" \" run and the long\n// \"the option to plot...\" run are untouched.\nconst phIntro = findParagraphByText(\n  items,\n  \"There are two options for plotting pH. The first is the option to plot the Seasonal Kendall trend line. It is important to note that even if a trend line will plot it does not mean the trend is statistically significant. Be sure to note the significance indicated in the sub-title of the plot. The second is to specify the applicable OWRD Basin specific pH criteria. Use this if your selected location is in one of the unique locations that has a separate standard from the rest of the basin (examples include Columbia or Snake River main-stem or Cascade Lakes > 3,000ft)\",\n  1\n);\nif (phIntro) {\n  phIntro.insertOoxml(\n    wrapParagraph(\n      '<w:p w:rsidR=\"00CC7C13\" w:rsidRDefault=\"00CC7C13\" w:rsidP=\"00CC7C13\">' +\n        \"<w:r><w:t>There are two options for plotting pH. The first is</w:t></w:r>\" +\n        '<w:r w:rsidRPr=\"00CC7C13\"><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        \"<w:r><w:t>the option to plot the Seasonal Kendall trend line. It is important to note that even if a trend line will plot it does not mean the trend is statistically significant. Be sure to note the significance indicated in the sub-title of the plot. The second is to specify the applicable OWRD Basin specific pH criteria. Use this if your selected location is in one of the unique locations that has a separate standard from the rest of the basin (examples include Columbia or Snake River main-stem or Cascade Lakes &gt; 3,000ft)</w:t></w:r>\" +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Intent (per commit message \"Modified handling of empty SeaKen results\"\n# and the supplied OOXML diff): remove the stray w:proofErr\n# (spellStart/spellEnd/gramStart/gramEnd) markers that Word's proofing\n# pass had inserted around certain runs, and merge the runs that had\n# been artificially split by those markers back into single runs with\n# the same visible text. No visible text content changes.\n\n$d = $word.ActiveDocument\n\n$W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\nfunction Get-ParagraphPlainText($range) {\n    # Paragraph.Range.Text includes the trailing paragraph mark (CR, 13)\n    # and, for the sole paragraph inside a table cell, a trailing cell\n    # mark (BEL, 7) as well. Strip those so we can match on visible text.\n    return $range.Text.TrimEnd([char]7, [char]13)\n}\n\nfunction Replace-ParagraphByText($targetText, $newInnerXml, $occurrence) {\n    $seen = 0\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = Get-ParagraphPlainText $p.Range\n        if ($t -eq $targetText) {\n            $seen = $seen + 1\n            if ($seen -eq $occurrence) {\n                $p.Range.InsertXML($newInnerXml)\n                return $true\n            }\n        }\n    }\n    return $false\n}\n\n# 1) First table: \"Station_ID\" cell - drop spellStart/spellEnd around the run.\nReplace-ParagraphByText \"Station_ID\" ('<w:p ' + $W_NS + ' w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\"><w:r><w:t>Station_ID</w:t></w:r></w:p>') 1 | Out-Null\n\n# 2) First table: \"Station_Description\" cell.\nReplace-ParagraphByText \"Station_Description\" ('<w:p ' + $W_NS + ' w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\"><w:r><w:t>Station_Description</w:t></w:r></w:p>') 1 | Out-Null\n\n# 3) First table: \"SampleType\" cell.\nReplace-ParagraphByText \"SampleType\" ('<w:p ' + $W_NS + ' w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\"><w:r><w:t>SampleType</w:t></w:r></w:p>') 1 | Out-Null\n\n# 4) First table: \"StatusIdentifier\" cell.\nReplace-ParagraphByText \"StatusIdentifier\" ('<w:p ' + $W_NS + ' w:rsidR=\"00825AC1\" w:rsidRDefault=\"00825AC1\" w:rsidP=\"00825AC1\"><w:r><w:t>StatusIdentifier</w:t></w:r></w:p>') 1 | Out-Null\n\n# 5) \"Seasonal Kendall Results: ...\" paragraph - merge the 3 runs that\n# spelled out \"...seasonal \", \"kendall\", \" trend analysis...\" (split by\n# spellStart/gramStart/spellEnd/gramEnd around \"kendall\") into one run.\n$seaKenText = \"Seasonal Kendall Results: Table detailing results of seasonal kendall trend analysis for each parameter at each station. Column descriptions follow:\"\n$seaKenXml = '<w:p ' + $W_NS + ' w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\"><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Seasonal Kendall Results: </w:t></w:r><w:r><w:t>Table detailing results of seasonal kendall trend analysis for each parameter at each station. Column descriptions follow:</w:t></w:r></w:p>'\nReplace-ParagraphByText $seaKenText $seaKenXml 1 | Out-Null\n\n# 6) Second table: \"Station_ID\" cell (2nd occurrence of that exact text).\nReplace-ParagraphByText \"Station_ID\" ('<w:p ' + $W_NS + ' w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\"><w:r><w:t>Station_ID</w:t></w:r></w:p>') 2 | Out-Null\n\n# 7) Second table: \"Pvalue\" cell.\nReplace-ParagraphByText \"Pvalue\" ('<w:p ' + $W_NS + ' w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\"><w:r><w:t>Pvalue</w:t></w:r></w:p>') 1 | Out-Null\n\n# 8) Second table: \"signif\" cell.\nReplace-ParagraphByText \"signif\" ('<w:p ' + $W_NS + ' w:rsidR=\"008F7447\" w:rsidRDefault=\"008F7447\" w:rsidP=\"00825AC1\"><w:r><w:t>signif</w:t></w:r></w:p>') 1 | Out-Null\n\n# 9) \"E. Coli/Enterococcus: \" heading - merge 3 runs into 1.\n$eColiXml = '<w:p ' + $W_NS + ' w:rsidR=\"00CC7C13\" w:rsidRDefault=\"00CC7C13\" w:rsidP=\"00CC7C13\"><w:pPr><w:pStyle w:val=\"Heading4\"/></w:pPr><w:r><w:t xml:space=\"preserve\">E. Coli/Enterococcus: </w:t></w:r></w:p>'\nReplace-ParagraphByText \"E. Coli/Enterococcus: \" $eColiXml 1 | Out-Null\n\n# 10) \"pH:\" heading - merge 2 runs into 1 (keep lastRenderedPageBreak).\n$phHeadingXml = '<w:p ' + $W_NS + ' w:rsidR=\"00CC7C13\" w:rsidRDefault=\"00CC7C13\" w:rsidP=\"00CC7C13\"><w:pPr><w:pStyle w:val=\"Heading4\"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>pH:</w:t></w:r></w:p>'\nReplace-ParagraphByText \"pH:\" $phHeadingXml 1 | Out-Null\n\n# 11) \"There are two options for plotting pH. The first is ...\" paragraph -\n# merge the first 3 runs (\"There are two options for plotting \", \"pH.\",\n# \" The first is\") into a single run; the trailing \" \" run and the long\n# \"the option to plot...\" run are untouched.\n$phIntroText = \"There are two options for plotting pH. The first is the option to plot the Seasonal Kendall trend line. It is important to note that even if a trend line will plot it does not mean the trend is statistically significant. Be sure to note the significance indicated in the sub-title of the plot. The second is to specify the applicable OWRD Basin specific pH criteria. Use this if your selected location is in one of the unique locations that has a separate standard from the rest of the basin (examples include Columbia or Snake River main-stem or Cascade Lakes > 3,000ft)\"\n$phIntroXml = '<w:p ' + $W_NS + ' w:rsidR=\"00CC7C13\" w:rsidRDefault=\"00CC7C13\" w:rsidP=\"00CC7C13\"><w:r><w:t>There are two options for plotting pH. The first is</w:t></w:r><w:r w:rsidRPr=\"00CC7C13\"><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>the option to plot the Seasonal Kendall trend line. It is important to note that even if a trend line will plot it does not mean the trend is statistically significant. Be sure to note the significance indicated in the sub-title of the plot. The second is to specify the applicable OWRD Basin specific pH criteria. Use this if your selected location is in one of the unique locations that has a separate standard from the rest of the basin (examples include Columbia or Snake River main-stem or Cascade Lakes &gt; 3,000ft)</w:t></w:r></w:p>'\nReplace-ParagraphByText $phIntroText $phIntroXml 1 | Out-Null\n"}
